$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written/stay as text, not auto-converted to numbers
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.397.99'
$ws.Range('E2').Value = '  -3.24%  '
$ws.Range('D3').Value = '3.060.78'
$ws.Range('E3').Value = '  -2.73%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '535.37'
$ws.Range('E5').Value = '  -4.65%  '
$ws.Range('D6').Value = '131.97'
$ws.Range('E6').Value = '  -9.44%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.050.84'
$ws.Range('E8').Value = '  -2.67%  '
$ws.Range('E9').Value = '  -2.49%  '
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  -2.67%  '
$ws.Range('D11').Value = '6.04'
$ws.Range('E11').Value = '  -10.57%  '
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').Value = '34.03'
$ws.Range('E14').Value = '  -8.45%  '
$ws.Range('D15').Value = '3.514.48'
$ws.Range('E15').Value = '  -3.89%  '
$ws.Range('D16').Value = '62.477.25'
$ws.Range('E16').Value = '  -3.38%  '
$ws.Range('E17').Value = '  -2.57%  '
$ws.Range('D18').Value = '3.063.67'
$ws.Range('E18').Value = '  -3.23%  '
$ws.Range('E19').Value = '  -5.44%  '
$ws.Range('D20').Value = '472.98'
$ws.Range('E20').Value = '  -8.53%  '
$ws.Range('D21').Value = '13.18'
$ws.Range('E21').Value = '  -6.72%  '
$ws.Range('D22').Value = '0.692'
$ws.Range('E22').Value = '  -4.00%  '
$ws.Range('E23').Value = '  -5.49%  '
$ws.Range('D24').Value = '78.12'
$ws.Range('E24').Value = '  -0.88%  '
$ws.Range('E25').Value = '  -8.32%  '
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.30%  '
$ws.Range('E27').Value = '  -6.18%  '
$ws.Range('D28').Value = '8.08'
$ws.Range('E28').Value = '  -10.85%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('E30').Value = '  -4.10%  '
$ws.Range('E31').Value = '  -16.32%  '
$ws.Range('D32').Value = '1.08'
$ws.Range('E32').Value = '  -4.94%  '
$ws.Range('E33').Value = '  -9.57%  '
$ws.Range('D34').Value = '55.91'
$ws.Range('E34').Value = '  +3.23%  '
$ws.Range('D35').Value = '5.86'
$ws.Range('E35').Value = '  -4.11%  '
$ws.Range('E36').Value = '  -5.04%  '
$ws.Range('D37').Value = '470.20'
$ws.Range('E37').Value = '  -13.92%  '
$ws.Range('D38').Value = '3.098.74'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('E39').Value = '  -11.03%  '
$ws.Range('D40').Value = '0.0782'
$ws.Range('E40').Value = '  -5.76%  '
$ws.Range('D41').Value = '7.96'
$ws.Range('E41').Value = '  -4.59%  '
$ws.Range('D42').Value = '0.111'
$ws.Range('E42').Value = '  -10.49%  '
$ws.Range('D43').Value = '2.56'
$ws.Range('E43').Value = '  -8.25%  '
$ws.Range('E45').Value = '  -9.26%  '
$ws.Range('D46').Value = '2.00'
$ws.Range('E46').Value = '  -9.97%  '
$ws.Range('D47').Value = '24.00'
$ws.Range('E47').Value = '  -6.70%  '
$ws.Range('D48').Value = '116.83'
$ws.Range('E48').Value = '  -5.01%  '
$ws.Range('E49').Value = '  -2.83%  '
$ws.Range('D50').Value = '0.0₃0508'
$ws.Range('E50').Value = '  -3.35%  '
$ws.Range('D51').Value = '1.96'
$ws.Range('E51').Value = '  -7.56%  '

# Restore default (Normal) style on the Price column so no residual number formatting remains
$ws.Range("D2:D51").Style = "Normal"

